$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.354.19'
$ws.Range('E2').Value = '  -7.61%  '
$ws.Range('D3').Value = '2.470.66'
$ws.Range('E3').Value = '  -12.51%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '465.05'
$ws.Range('E5').Value = '  -7.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.64'
$ws.Range('E6').Value = '  -2.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.490'
$ws.Range('E8').Value = '  -6.76%  '
$ws.Range('D9').Value = '2.469.40'
$ws.Range('E9').Value = '  -12.46%  '
$ws.Range('E10').Value = '  -6.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.35'
$ws.Range('E11').Value = '  -9.47%  '
$ws.Range('E12').Value = '  -7.64%  '
$ws.Range('E13').Value = '  -4.07%  '
$ws.Range('D14').Value = '2.877.36'
$ws.Range('E14').Value = '  -13.35%  '
$ws.Range('D15').Value = '54.205.48'
$ws.Range('E15').Value = '  -8.10%  '
$ws.Range('E16').Value = '  +1.56%  '
$ws.Range('E17').Value = '  -7.57%  '
$ws.Range('D18').Value = '2.472.34'
$ws.Range('E18').Value = '  -12.56%  '
$ws.Range('E19').Value = '  -10.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '312.08'
$ws.Range('E20').Value = '  -9.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.46'
$ws.Range('E21').Value = '  -13.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.994'
$ws.Range('E22').Value = '  -0.66%  '
$ws.Range('E23').Value = '  +1.11%  '
$ws.Range('E24').Value = '  -13.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '56.65'
$ws.Range('E25').Value = '  -9.78%  '
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('D27').Value = '2.558.61'
$ws.Range('E27').Value = '  -13.20%  '
$ws.Range('E28').Value = '  -9.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.155'
$ws.Range('E29').Value = '  -8.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.21'
$ws.Range('E30').Value = '  -1.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.996'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').Value = '0.0₃0728'
$ws.Range('E32').Value = '  -7.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '150.97'
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('E34').Value = '  -6.56%  '
$ws.Range('E35').Value = '  -10.46%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.02'
$ws.Range('E36').Value = '  -4.80%  '
$ws.Range('E37').Value = '  -13.15%  '
$ws.Range('E38').Value = '  -4.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.791'
$ws.Range('E39').Value = '  -11.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '33.58'
$ws.Range('E40').Value = '  -8.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('E42').Value = '  -3.20%  '
$ws.Range('E43').Value = '  -4.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.28'
$ws.Range('E44').Value = '  -5.94%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.16'
$ws.Range('E45').Value = '  -1.85%  '
$ws.Range('E46').Value = '  -7.15%  '
$ws.Range('D47').Value = '1.963.50'
$ws.Range('E47').Value = '  -11.82%  '
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0872'
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.33'
$ws.Range('E50').Value = '  -4.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.76'
$ws.Range('E51').Value = '  -12.16%  '
